# Auto-generated edit script applying the Garuda_Profits.xlsx diff
# Updates currentAveragePrice / LevePrice / LeveProfit columns (H-N) for affected leves
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H16").Value = 46706.668
$ws.Range("I16").Value = 100
$ws.Range("J16").Value = 70010
$ws.Range("K16").Value = 100
$ws.Range("L16").Value = 70010
$ws.Range("M16").Value = 130
$ws.Range("N16").Value = -70470

$ws.Range("H62").Value = 1995
$ws.Range("I62").Value = 1995
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 1995
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -1371
$ws.Range("N62").ClearContents()

$ws.Range("H64").Value = 3028.25
$ws.Range("I64").Value = 3006.125
$ws.Range("J64").Value = 3072.5
$ws.Range("K64").Value = 3006.125
$ws.Range("L64").Value = 3072.5
$ws.Range("M64").Value = -2758.125
$ws.Range("N64").Value = -3568.5

$ws.Range("H65").Value = 1995
$ws.Range("I65").Value = 1995
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9975
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -6855
$ws.Range("N65").ClearContents()

$ws.Range("H67").Value = 3028.25
$ws.Range("I67").Value = 3006.125
$ws.Range("J67").Value = 3072.5
$ws.Range("K67").Value = 3006.125
$ws.Range("L67").Value = 3072.5
$ws.Range("M67").Value = -2148.125
$ws.Range("N67").Value = -4788.5

$ws.Range("H69").Value = 3428.261
$ws.Range("I69").Value = 3390.9092
$ws.Range("J69").Value = 4250
$ws.Range("K69").Value = 10172.7276
$ws.Range("L69").Value = 12750
$ws.Range("M69").Value = -9298.7276

$ws.Range("H72").Value = 3428.261
$ws.Range("I72").Value = 3390.9092
$ws.Range("J72").Value = 4250
$ws.Range("K72").Value = 30518.1828
$ws.Range("L72").Value = 38250
$ws.Range("M72").Value = -26150.1828

$ws.Range("H74").Value = 3824.825
$ws.Range("I74").Value = 3776.1904
$ws.Range("J74").Value = 3878.5789
$ws.Range("K74").Value = 3776.1904
$ws.Range("L74").Value = 3878.5789
$ws.Range("M74").Value = -2840.1904
$ws.Range("N74").Value = -5750.5789

$ws.Range("H77").Value = 3824.825
$ws.Range("I77").Value = 3776.1904
$ws.Range("J77").Value = 3878.5789
$ws.Range("K77").Value = 18880.952
$ws.Range("L77").Value = 19392.8945
$ws.Range("M77").Value = -14200.952
$ws.Range("N77").Value = -28752.8945

$ws.Range("H94").Value = 5000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 5000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 5000
$ws.Range("N94").Value = -5902

$ws.Range("H137").Value = 1993.1177
$ws.Range("I137").Value = 1740
$ws.Range("J137").Value = 2600.6
$ws.Range("K137").Value = 5220
$ws.Range("L137").Value = 7801.799999999999
$ws.Range("M137").Value = -2670
$ws.Range("N137").Value = -12901.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 311.85715
$ws.Range("I22").Value = 312.5
$ws.Range("J22").Value = 299
$ws.Range("K22").Value = 312.5
$ws.Range("L22").Value = 299
$ws.Range("M22").Value = -139.5

$ws.Range("H105").Value = 2738.889
$ws.Range("I105").Value = 2738.889
$ws.Range("J105").Value = 0
$ws.Range("K105").Value = 2738.889
$ws.Range("L105").Value = 0
$ws.Range("M105").Value = -991.8890000000001
$ws.Range("N105").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 465.47058
$ws.Range("I22").Value = 485.9
$ws.Range("J22").Value = 436.2857
$ws.Range("K22").Value = 485.9
$ws.Range("L22").Value = 436.2857
$ws.Range("M22").Value = -135.9
$ws.Range("N22").Value = -1136.2857

$ws.Range("H31").Value = 9092665
$ws.Range("I31").Value = 1857.6666
$ws.Range("J31").Value = 50001300
$ws.Range("K31").Value = 1857.6666
$ws.Range("L31").Value = 50001300
$ws.Range("M31").Value = -1562.6666

$ws.Range("H34").Value = 9092665
$ws.Range("I34").Value = 1857.6666
$ws.Range("J34").Value = 50001300
$ws.Range("K34").Value = 1857.6666
$ws.Range("L34").Value = 50001300
$ws.Range("M34").Value = -1655.6666

$ws.Range("H62").Value = 111113080
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 111113080
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 111113080
$ws.Range("N62").Value = -111114328

$ws.Range("H65").Value = 111113080
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 111113080
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 555565400
$ws.Range("N65").Value = -555571640

$ws.Range("H107").Value = 614.5909
$ws.Range("I107").Value = 557.2857
$ws.Range("J107").Value = 714.875
$ws.Range("K107").Value = 557.2857
$ws.Range("L107").Value = 714.875
$ws.Range("M107").Value = 1362.7143
$ws.Range("N107").Value = -4554.875

$ws.Range("H135").Value = 32890
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 32890
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 32890
$ws.Range("N135").Value = -43030

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 712.3125
$ws.Range("I5").Value = 317.14285
$ws.Range("J5").Value = 1019.6667
$ws.Range("K5").Value = 951.4285500000001
$ws.Range("L5").Value = 3059.0001
$ws.Range("M5").Value = -839.4285500000001
$ws.Range("N5").Value = -3283.0001

$ws.Range("H135").Value = 712.3125
$ws.Range("I135").Value = 317.14285
$ws.Range("J135").Value = 1019.6667
$ws.Range("K135").Value = 2854.28565
$ws.Range("L135").Value = 9177.0003
$ws.Range("M135").Value = -319.2856500000003
$ws.Range("N135").Value = -14247.0003

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1725.1428
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1725.1428
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1725.1428
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -2101.1428

$ws.Range("H68").Value = 1347.4546
$ws.Range("I68").Value = 1190.25
$ws.Range("J68").Value = 1766.6666
$ws.Range("K68").Value = 1190.25
$ws.Range("L68").Value = 1766.6666
$ws.Range("M68").Value = -441.25
$ws.Range("N68").Value = -3264.6666

$ws.Range("H69").Value = 48003
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 48003
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 48003
$ws.Range("N69").Value = -49625

$ws.Range("H71").Value = 1347.4546
$ws.Range("I71").Value = 1190.25
$ws.Range("J71").Value = 1766.6666
$ws.Range("K71").Value = 5951.25
$ws.Range("L71").Value = 8833.333000000001
$ws.Range("M71").Value = -2207.25
$ws.Range("N71").Value = -16321.333

$ws.Range("H72").Value = 48003
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 48003
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 144009
$ws.Range("N72").Value = -152121

$ws.Range("H82").Value = 2097.6
$ws.Range("I82").Value = 2629.3333
$ws.Range("J82").Value = 1300
$ws.Range("K82").Value = 2629.3333
$ws.Range("L82").Value = 1300
$ws.Range("M82").Value = -2268.3333
$ws.Range("N82").Value = -2022

$ws.Range("H85").Value = 2097.6
$ws.Range("I85").Value = 2629.3333
$ws.Range("J85").Value = 1300
$ws.Range("K85").Value = 2629.3333
$ws.Range("L85").Value = 1300
$ws.Range("M85").Value = -1381.3333
$ws.Range("N85").Value = -3796

$ws.Range("H93").Value = 1082604.4
$ws.Range("I93").Value = 1803414.1
$ws.Range("J93").Value = 1389.7
$ws.Range("K93").Value = 1803414.1
$ws.Range("L93").Value = 1389.7
$ws.Range("M93").Value = -1802166.1
$ws.Range("N93").Value = -3885.7

$ws.Range("H100").Value = 1166
$ws.Range("I100").Value = 937.5
$ws.Range("J100").Value = 2994
$ws.Range("K100").Value = 937.5
$ws.Range("L100").Value = 2994
$ws.Range("M100").Value = -396.5
$ws.Range("N100").Value = -4076

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H26").Value = 6483
$ws.Range("I26").Value = 4310.6665
$ws.Range("J26").Value = 13000
$ws.Range("K26").Value = 4310.6665
$ws.Range("L26").Value = 13000
$ws.Range("M26").Value = -4017.6665
$ws.Range("N26").Value = -13586

$ws.Range("H63").Value = 14666.667
$ws.Range("I63").Value = 0
$ws.Range("J63").Value = 14666.667
$ws.Range("K63").Value = 0
$ws.Range("L63").Value = 14666.667
$ws.Range("N63").Value = -15914.667

$ws.Range("H66").Value = 14666.667
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 14666.667
$ws.Range("K66").Value = 0
$ws.Range("L66").Value = 44000.001
$ws.Range("N66").Value = -50240.001

$ws.Range("H69").Value = 20000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 20000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 20000
$ws.Range("N69").Value = -21498

$ws.Range("H72").Value = 20000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 20000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 60000
$ws.Range("N72").Value = -67488

$ws.Range("H76").Value = 19250
$ws.Range("I76").Value = 19000
$ws.Range("J76").Value = 19500
$ws.Range("K76").Value = 19000
$ws.Range("L76").Value = 19500
$ws.Range("M76").Value = -18685
$ws.Range("N76").Value = -20130

$ws.Range("H79").Value = 19250
$ws.Range("I79").Value = 19000
$ws.Range("J79").Value = 19500
$ws.Range("K79").Value = 19000
$ws.Range("L79").Value = 19500
$ws.Range("M79").Value = -17908
$ws.Range("N79").Value = -21684

$ws.Range("H80").Value = 0
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 0
$ws.Range("N80").ClearContents()

$ws.Range("H82").Value = 20000
$ws.Range("I82").Value = 0
$ws.Range("J82").Value = 20000
$ws.Range("K82").Value = 0
$ws.Range("L82").Value = 20000
$ws.Range("N82").Value = -20766

$ws.Range("H83").Value = 0
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 0
$ws.Range("N83").ClearContents()

$ws.Range("H85").Value = 20000
$ws.Range("I85").Value = 0
$ws.Range("J85").Value = 20000
$ws.Range("K85").Value = 0
$ws.Range("L85").Value = 20000
$ws.Range("N85").Value = -22652

$ws.Range("H87").Value = 61818.184
$ws.Range("I87").Value = 0
$ws.Range("J87").Value = 61818.184
$ws.Range("K87").Value = 0
$ws.Range("L87").Value = 61818.184
$ws.Range("N87").Value = -64314.184

$ws.Range("H90").Value = 61818.184
$ws.Range("I90").Value = 0
$ws.Range("J90").Value = 61818.184
$ws.Range("K90").Value = 0
$ws.Range("L90").Value = 185454.552
$ws.Range("N90").Value = -197934.552
